$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParagraphIndexByText($target) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $target) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXml($oldText, $newParagraphInnerXml) {
    $idx = Find-ParagraphIndexByText $oldText
    if ($idx -eq -1) {
        throw "Could not find paragraph with text: $oldText"
    }
    $p = $d.Paragraphs.Item($idx)
    $xml = "<w:p $wns>$newParagraphInnerXml</w:p>"
    [void]$p.Range.InsertXML($xml)
}

function Remove-ParagraphWithText($target) {
    $idx = Find-ParagraphIndexByText $target
    if ($idx -eq -1) {
        throw "Could not find paragraph to remove with text: $target"
    }
    [void]$d.Paragraphs.Item($idx).Range.Delete()
}

$euro = [char]8364

# 1a. H1 title
Set-ParagraphXml `
    "Play Master Joker Free | Review of Pragmatic Play Slot Game" `
    "<w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Play Master Joker Free - Exciting Slot Game with Impressive Graphics</w:t></w:r>"

# 2. "What we like" bullet list items
Set-ParagraphXml `
    "Stylish and modern design with refreshing fruit theme" `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Simple gameplay with one payline and high RTP</w:t></w:r>"

Set-ParagraphXml `
    ("High RTP with a chance to win up to " + $euro + "1,000,000") `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Impressive graphics and vibrant colors</w:t></w:r>"

Set-ParagraphXml `
    "Multiplayer Wheel hidden feature with 2x to 100x multiplier potential" `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Exciting Multiplayer Wheel and Wild Symbol features</w:t></w:r>"

Set-ParagraphXml `
    "Accessible on mobile, tablet and desktop devices" `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Chance to win up to 10,000 times the bet</w:t></w:r>"

# 3. "What we don't like" bullet list: first item text changes, second item is removed entirely
Set-ParagraphXml `
    "Limited special features compared to other online slots" `
    "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r/><w:r><w:t>Master Joker has few features compared to other slots</w:t></w:r>"

Remove-ParagraphWithText "Only one payline to bet on"

# 4a. Bold "title" paragraph near the end
Set-ParagraphXml `
    "Play Master Joker Free | Review of Pragmatic Play Slot Game" `
    "<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Master Joker Free - Exciting Slot Game with Impressive Graphics</w:t></w:r>"

# 4b. Italic summary/description paragraph at the end
Set-ParagraphXml `
    "Read our review of Master Joker, an online slot game by Pragmatic Play. Play it for free or real money and activate the Multiplayer Wheel feature." `
    "<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Master Joker, a visually stunning slot game with multiplayer and Wild Symbol features.</w:t></w:r>"

Write-Output "Done"
